$d = $word.ActiveDocument

# Mapping of old division expressions to new ones, applied via Find/Replace
# on the whole document content. Each "old" text is unique in the document,
# so replacing all occurrences safely targets exactly the intended cell.
$replacements = @(
    @("12÷4=", "20÷7="),
    @("59÷2=", "94÷6="),
    @("58÷6=", "35÷3="),
    @("51÷5=", "73÷9="),
    @("81÷5=", "92÷5="),
    @("52÷4=", "22÷4="),
    @("22÷8=", "92÷4="),
    @("34÷8=", "40÷4="),
    @("67÷2=", "92÷4="),
    @("24÷5=", "88÷3="),
    @("16÷2=", "21÷2="),
    @("55÷8=", "99÷8="),
    @("74÷3=", "49÷2="),
    @("26÷5=", "27÷7="),
    @("32÷4=", "54÷5="),
    @("59÷6=", "89÷4="),
    @("30÷8=", "25÷3="),
    @("99÷7=", "91÷7="),
    @("89÷7=", "64÷3="),
    @("18÷7=", "10÷9="),
    @("22÷9=", "95÷4="),
    @("56÷7=", "56÷2="),
    @("15÷9=", "20÷9="),
    @("50÷7=", "82÷2="),
    @("27÷5=", "65÷7="),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

